$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add headers for columns I (I0) and J (IF), copying style/format from H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Populate I/J values per row (I0 and IF columns)
$data = @{
    2 = @(9, 9)
    3 = @(9, 9)
    4 = @(10, 10)
    5 = @(9, 9)
    6 = @(9, 9)
    7 = @(9, 9)
    8 = @(9, 9)
    9 = @(9, 9)
    10 = @(8, 9)
    11 = @(9, 9)
    12 = @(8, 9)
    13 = @(8, 9)
    14 = @(9, 9)
    15 = @(9, 9)
    16 = @(9, 10)
    17 = @(9, 9)
    18 = @(9, 9)
    19 = @(8, 9)
    20 = @(9, 9)
    21 = @(10, 10)
    22 = @(9, 9)
    23 = @(9, 9)
    24 = @(9, 9)
    25 = @(9, 9)
    26 = @(9, 9)
    27 = @(8, 8)
    28 = @(9, 9)
    29 = @(9, 9)
    30 = @(9, 9)
    31 = @(9, 9)
    32 = @(9, 9)
    33 = @(9, 9)
    34 = @(9, 9)
    35 = @(9, 9)
    36 = @(9, 9)
    37 = @(9, 9)
    38 = @(9, 9)
    39 = @(9, 10)
    40 = @(9, 9)
    41 = @(9, 9)
    42 = @(9, 9)
    43 = @(9, 9)
    44 = @(9, 9)
    45 = @(9, 9)
    46 = @(9, 9)
    47 = @(9, 10)
    48 = @(9, 9)
    49 = @(10, 10)
    50 = @(9, 9)
    51 = @(9, 9)
    52 = @(9, 9)
    53 = @(9, 9)
    54 = @(9, 9)
    55 = @(9, 9)
    56 = @(10, 10)
    57 = @(9, 9)
    58 = @(9, 9)
    59 = @(9, 9)
    60 = @(10, 10)
    61 = @(8, 9)
    62 = @(9, 9)
    63 = @(9, 9)
    64 = @(9, 9)
    65 = @(9, 9)
    66 = @(8, 8)
    67 = @(7, 7)
    68 = @(9, 9)
    69 = @(7, 7)
    70 = @(8, 8)
    71 = @(5, 5)
    72 = @(6, 6)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 9).Value = $vals[0]
    $ws.Cells.Item($r, 10).Value = $vals[1]
}
